$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename header row (row 1) texts:
#    "<Column>_old" -> "<Column>_FV2404"
#    "<Column>_new" -> "<Column>_FV2410"
$headerRange = $ws.Range("A1:U1")
for ($i = 1; $i -le $headerRange.Columns.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $text = $cell.Value()
    if ($text -ne $null) {
        $newText = $text -replace "_old$", "_FV2404"
        $newText = $newText -replace "_new$", "_FV2410"
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}

# 2) Turn the used range into an Excel Table ("Table1") with autofilter,
#    now that the header cells carry the final (FV2404/FV2410) names.
$tableRange = $ws.Range("A1:U64")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# 3) Freeze the header row (split beneath row 1 / above row 2).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
